$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = 1
$ws.Range("D7").Value = ""

$ws.Range("D10").Value = ""

$ws.Range("C12").Value = 0
$ws.Range("D12").Value = ""

$ws.Range("C14").Value = 0
$ws.Range("D14").Value = ""

$ws.Range("C15").Value = 0
$ws.Range("D15").Value = ""

$ws.Range("C16").Value = 1
$ws.Range("D16").Value = ""

$ws.Range("C19").Value = 0
$ws.Range("D19").Value = "LED?"

$ws.Range("C20").Value = 0
$ws.Range("D20").Value = ""

$ws.Range("C22").Value = 1
$ws.Range("D22").Value = ""

$ws.Range("C23").Select()
